# Add new rows for "Ukenglish" and "Usenglish" to the Google summary sheet,
# mirroring the layout/formatting of the existing language rows
# (Category, edit distance mean, edit distance stdev, accuracy mean, accuracy stdev).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Ukenglish
$ws.Range("A17").Value = "Ukenglish"
$ws.Range("B17").Value = 42.5
$ws.Range("C17").Value = 32
$ws.Range("D17").Value = 13
$ws.Range("F17").Value = 9.1300000000000008
$ws.Range("D17").Interior.Color = $ws.Range("D15").Interior.Color

# Row 18: Usenglish
$ws.Range("A18").Value = "Usenglish"
$ws.Range("B18").Value = 17.75
$ws.Range("C18").Value = 6.71
$ws.Range("D18").Value = 4.625
$ws.Range("F18").Value = 2.5590000000000002
$ws.Range("D18").Interior.Color = $ws.Range("D15").Interior.Color

# Match the final selection/view state from the diff.
$ws.Range("C18").Select()
